# Update WS_holdings.xlsx:
#  - Refresh the "as of" date in the confidentiality notice (A16)
#  - Refresh the Weight (D) / Percent Change (E) values for rows 2-13
#
# The worksheet ships sheet-protected, so editing the locked D:E cells
# requires a temporary Unprotect/Protect cycle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Confidentiality footer: bump the "as of" date ----------------------
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

# --- Weight (D) / Percent Change (E) refresh, rows 2-13 ------------------
$ws.Range("D2").Value = 0.03092532749574342
$ws.Range("E2").Value = 0.003780068728522323

$ws.Range("D3").Value = 0.02422203379731474
$ws.Range("E3").Value = -0.001390176088971073

$ws.Range("D4").Value = 0.05176311621840524
$ws.Range("E4").Value = -0.0007058823529412228

$ws.Range("D5").Value = 0.1397636377900123
$ws.Range("E5").Value = 0.006135959954787662

$ws.Range("D6").Value = 0.02835806798318123
$ws.Range("E6").Value = 0.01102362204724394

$ws.Range("D7").Value = 0.1209843135667486
$ws.Range("E7").Value = 0.0119077463023316

$ws.Range("D8").Value = 0.1006709622291646
$ws.Range("E8").Value = 0.01236881559220393

$ws.Range("D9").Value = 0.02774503029306521
$ws.Range("E9").Value = 0.01799824407374895

$ws.Range("D10").Value = 0.1191924858343993
$ws.Range("E10").Value = 0.0206371468643558

$ws.Range("D11").Value = 0.2541670502630008
$ws.Range("E11").Value = 0.01508524502365427

$ws.Range("D12").Value = 0.1022079745289646
$ws.Range("E12").Value = 0.01256133464180564

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.01197990755238787

$ws.Protect()
